$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Сочи"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2024-11-25"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "09:30"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = "В"
$ws.Range("G2").Value = 1016
$ws.Range("J2").Value = 2178

# Row 3
$ws.Range("A3").Value = "Сочи"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2024-11-25"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "09:30"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = "В"
$ws.Range("G3").Value = 1016
$ws.Range("J3").Value = 2177

# Row 4
$ws.Range("A4").Value = "Сочи"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2024-11-25"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = "09:30"
$ws.Range("D4").Value = 10
$ws.Range("G4").Value = 1016
$ws.Range("J4").Value = 2176

# Row 5
$ws.Range("A5").Value = "Сочи"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2024-11-25"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = "09:30"
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = "В"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 1016
$ws.Range("J5").Value = 2175

# Row 6
$ws.Range("A6").Value = "Владивосток"
$ws.Range("C6").Value = "15:15"
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = "С"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1011
$ws.Range("J6").Value = 2174

# Row 7
$ws.Range("A7").Value = "Владивосток"
$ws.Range("C7").Value = "15:15"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "С"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1011
$ws.Range("J7").Value = 2173

# Row 8
$ws.Range("A8").Value = "Владивосток"
$ws.Range("C8").Value = "15:15"
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = "С"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1011
$ws.Range("J8").Value = 2172

# Row 9
$ws.Range("C9").Value = "15:15"
$ws.Range("G9").Value = 990
$ws.Range("J9").Value = 2171

# Row 10
$ws.Range("C10").Value = "15:15"
$ws.Range("G10").Value = 990
$ws.Range("J10").Value = 2170

# Row 11
$ws.Range("C11").Value = "15:15"
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 990
$ws.Range("J11").Value = 2169

